$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers for columns I and J
$ws.Range("I1").Value = "Course Name"
$ws.Range("J1").Value = "Student ID"

# Fill in data for rows 2-7
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 9).Value = "Associate Degree of Information Technology"
    $ws.Cells.Item($r, 10).Value = "A00123456"
}

# Copy the formatting already used in column H (style index 1) onto the
# newly added I:J columns so the new cells match the sheet's existing look.
$ws.Range("H1:H7").Copy()
$ws.Range("I1:J7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Select I2 to match the saved selection state
$ws.Range("I2").Select()
